$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @("ECs", "Il1b", "Il1r1", "ECs", 3, 1, 1802.141101666667, 5406.423305, 0.8789810559109003, 0.8789810559109003, 3, 1, 6.467363, 19.402089, 0.1496800897841894, 0.1496800897841894, 11655.10068169824, 104895.9061352842, 0.1315659633673452, 0.1315659633673452),
    @("ECs", "Il1b", "Il1r1", "FAPs", 3, 1, 1802.141101666667, 5406.423305, 0.8789810559109003, 0.8789810559109003, 3, 1, 26.306328, 78.91898400000001, 0.6088313794868691, 0.6088313794868691, 47407.71492272469, 426669.4343045222, 0.5351512488130582, 0.5351512488130582),
    @("ECs", "Il1b", "Il1r1", "M2", 3, 1, 1802.141101666667, 5406.423305, 0.8789810559109003, 0.8789810559109003, 2, 0.6666666666666666, 0.294709, 0.884127, 0.006820719600895863, 0.006820719600895864, 531.1072019310817, 4779.964817379735, 0.00599528331686762, 0.00599528331686762),
    @("ECs", "Il1b", "Il1r1", "sCs", 3, 1, 1802.141101666667, 5406.423305, 0.8789810559109003, 0.8789810559109003, 3, 1, 10.13950433333333, 30.418513, 0.2346678111280456, 0.2346678111280456, 18272.81750962728, 164455.3575866455, 0.2062685604136293, 0.2062685604136293),
    @("M2", "Il1b", "Il1r1", "ECs", 3, 1, 248.1049756666667, 744.314927, 0.121011375461416, 0.121011375461416, 3, 1, 6.467363, 19.402089, 0.1496800897841894, 0.1496800897841894, 1604.5849397425, 14441.2644576825, 0.018112993543973, 0.018112993543973),
    @("M2", "Il1b", "Il1r1", "FAPs", 3, 1, 248.1049756666667, 744.314927, 0.121011375461416, 0.121011375461416, 3, 1, 26.306328, 78.91898400000001, 0.6088313794868691, 0.6088313794868691, 6526.730868319353, 58740.57781487417, 0.07367552265577736, 0.07367552265577736),
    @("M2", "Il1b", "Il1r1", "M2", 3, 1, 248.1049756666667, 744.314927, 0.121011375461416, 0.121011375461416, 2, 0.6666666666666666, 0.294709, 0.884127, 0.006820719600895863, 0.006820719600895864, 73.11876927374766, 658.068923463729, 0.0008253846605410487, 0.0008253846605410489),
    @("M2", "Il1b", "Il1r1", "sCs", 3, 1, 248.1049756666667, 744.314927, 0.121011375461416, 0.121011375461416, 3, 1, 10.13950433333333, 30.418513, 0.2346678111280456, 0.2346678111280456, 2515.661475893728, 22640.95328304355, 0.02839747460112458, 0.02839747460112459),
    @("sCs", "Il1b", "Il1r1", "ECs", 1, 0.3333333333333333, 0.01551766666666667, 0.046553, 0.000007568627683662319, 0.00000756862768366232, 3, 1, 6.467363, 19.402089, 0.1496800897841894, 0.1496800897841894, 0.1003583832463333, 0.903225449217, 0.000001132872871233677, 0.000001132872871233677),
    @("sCs", "Il1b", "Il1r1", "FAPs", 1, 0.3333333333333333, 0.01551766666666667, 0.046553, 0.000007568627683662319, 0.00000756862768366232, 3, 1, 26.306328, 78.91898400000001, 0.6088313794868691, 0.6088313794868691, 0.4082128291280001, 3.673915462152, 0.000004608018033466636, 0.000004608018033466637),
    @("sCs", "Il1b", "Il1r1", "M2", 1, 0.3333333333333333, 0.01551766666666667, 0.046553, 0.000007568627683662319, 0.00000756862768366232, 2, 0.6666666666666666, 0.294709, 0.884127, 0.006820719600895863, 0.006820719600895864, 0.004573196025666667, 0.041158764231, 0.00000005162348719383863, 0.00000005162348719383864),
    @("sCs", "Il1b", "Il1r1", "sCs", 1, 0.3333333333333333, 0.01551766666666667, 0.046553, 0.000007568627683662319, 0.00000756862768366232, 3, 1, 10.13950433333333, 30.418513, 0.2346678111280456, 0.2346678111280456, 0.1573414484098889, 1.416073035689, 0.000001776113291768167, 0.000001776113291768167)
)

$r = 2
foreach ($rowVals in $rowsData) {
    $c = 1
    foreach ($val in $rowVals) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
